$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old detail rows (R0001..R0003 rows + the extra R0001 row),
# collapsing the table down to header + a single data row.
$ws.Rows("3:6").Delete()

# Update the remaining data row to the new shipment values.
$ws.Range("A2").Value = "R0004"
$ws.Range("B2").Value = "10254-ARI-I"
$ws.Range("C2").Value = 7

# Widen column A so the Remessa codes aren't truncated.
$ws.Columns("A").ColumnWidth = 19.1666666666667

# Leave the selection where the author left it when saving.
[void]$ws.Range("E15").Select()
